$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing row 3 values
$ws.Range("B3").Value = 20
$ws.Range("D3").Value = 45689.68715491898

# Update existing row 4 values
$ws.Range("B4").Value = 9
$ws.Range("D4").Value = 45689.68215277778

# Add new row 5
$ws.Range("A5").Value = "MAT141"
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = 45689.68708708333
$ws.Range("D5").Value = 45689.68714495492

# Copy style from C4/D4 (style index 2, date format) to C5/D5
$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("D4").Copy()
$ws.Range("D5").PasteSpecial(-4122) # xlPasteFormats
